$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A2").Value = "Hans"
$ws.Range("B2").Value = "Wurst"
$ws.Range("C2").Value = 1234

$ws.Range("F13").Select()
